$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Columns.Item(1).Insert()

# Set header for new ID column, matching the header style used by the rest of row 1
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set ID values for each row
$ws.Range("A2").Value = "RM 2"
$ws.Range("A3").Value = "RM 8"
$ws.Range("A4").Value = "RM 9"
$ws.Range("A5").Value = "RM 14"
$ws.Range("A6").Value = "RM 21"
$ws.Range("A7").Value = "RM 32"
$ws.Range("A8").Value = "RM 38"
$ws.Range("A9").Value = "RM 42"
$ws.Range("A10").Value = "RM 52 a"
$ws.Range("A11").Value = "RM 58"
$ws.Range("A12").Value = "RM 81"
$ws.Range("A13").Value = "RM 88"
$ws.Range("A14").Value = "RM 90"
$ws.Range("A15").Value = "RM 95"
$ws.Range("A16").Value = "RM 103"
$ws.Range("A17").Value = "RM 116"
$ws.Range("A18").Value = "RM 120"
$ws.Range("A19").Value = "RM 125"
$ws.Range("A20").Value = "RM 134"
$ws.Range("A21").Value = "RM 135"
$ws.Range("A22").Value = "RM 138"
$ws.Range("A23").Value = "RM 140"
$ws.Range("A24").Value = "RM 142a"
$ws.Range("A25").Value = "RM 145"
$ws.Range("A26").Value = "RM 232"
$ws.Range("A27").Value = "SC 5"
$ws.Range("A28").Value = "SC 92"
$ws.Range("A29").Value = "SC 101"
$ws.Range("A30").Value = "SC 105"
$ws.Range("A31").Value = "SC 119"
$ws.Range("A32").Value = "SC 120"
$ws.Range("A33").Value = "SC 132"
$ws.Range("A34").Value = "SC 193"
$ws.Range("A35").Value = "SC 232"
